# Apply the OT hostel excel upload sample update:
# - Add new leading column "course_master_pk" (values = 11) and shift the
#   existing "user_name" / "hostel_room_name" columns to B / C.
# - Populate 8 data rows beneath the header.
# - Header row keeps its yellow highlight; data rows get a wrapped,
#   vertically centered style.
# - Resize the three columns and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----------------------------------------------------------
$ws.Range("A1").Value = "course_master_pk"
$ws.Range("B1").Value = "user_name"
$ws.Range("C1").Value = "hostel_room_name"

# ---- Data rows ------------------------------------------------------------
$data = @(
    @(11, "JayasreePradhan",     "SILV-102"),
    @(11, "mundrawhat",          "MAHA-104"),
    @(11, "Akankshapathak1509",  "MAHA-104"),
    @(11, "DeepeshKaira",        "SILV-302"),
    @(11, "prakharkr29",         "SILV-302"),
    @(11, "negiitushar",         "SILV-302"),
    @(11, "Rana_ananya04",       "SILV-302"),
    @(11, "suramyasharma",       "SILV-301")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---- Column widths ---------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(2).ColumnWidth = 14.63
$ws.Columns.Item(3).ColumnWidth = 21.18

# ---- Styling --------------------------------------------------------------
# Header: A1/B1 already carry the original yellow-highlight style; copy that
# same formatting onto the newly introduced C1 header cell so all three
# header cells share one style (instead of minting a near-duplicate fill).
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: wrap text, vertically centered.
$lastRow = 1 + $data.Count
$dataRange = $ws.Range("A2:C" + $lastRow)
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

# ---- Selection -------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
